$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 2.6
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 401
$ws.Range("AT8").Value = 2.63
$ws.Range("AU8").Value = 8.5

$ws.Range("G9").Value = 1.75
$ws.Range("I9").Value = 5.5
$ws.Range("J9").Value = 2.5
$ws.Range("U9").Value = 2.5
$ws.Range("V9").Value = 1.5
$ws.Range("W9").Value = 4.75
$ws.Range("X9").Value = 6.5
$ws.Range("AH9").Value = 10
$ws.Range("AK9").Value = 67
$ws.Range("AZ9").Value = 151

$ws.Range("G19").Value = 1.62
$ws.Range("H19").Value = 3.8
$ws.Range("I19").Value = 5.25
$ws.Range("J19").Value = 2.2
$ws.Range("K19").Value = 2.3
$ws.Range("L19").Value = 5
$ws.Range("U19").Value = 1.73
$ws.Range("V19").Value = 2
$ws.Range("AC19").Value = 13
$ws.Range("AI19").Value = 26
$ws.Range("AN19").Value = 3.75
$ws.Range("AP19").Value = 19
$ws.Range("AQ19").Value = 26
$ws.Range("AY19").Value = 29
$ws.Range("AZ19").Value = 81

$ws.Range("G22").Value = 1.73
$ws.Range("H22").Value = 3.7
$ws.Range("J22").Value = 2.4
$ws.Range("K22").Value = 2.1
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = 1.07
$ws.Range("N22").Value = 9
$ws.Range("O22").Value = 1.33
$ws.Range("P22").Value = 3.25
$ws.Range("Q22").Value = 2.08
$ws.Range("R22").Value = 1.73
$ws.Range("S22").Value = 1.44
$ws.Range("T22").Value = 2.63
$ws.Range("U22").Value = 2
$ws.Range("V22").Value = 1.73
$ws.Range("W22").Value = 6
$ws.Range("Y22").Value = 8.5
$ws.Range("Z22").Value = 13
$ws.Range("AA22").Value = 15
$ws.Range("AB22").Value = 29
$ws.Range("AC22").Value = 9
$ws.Range("AD22").Value = 7
$ws.Range("AE22").Value = 19
$ws.Range("AF22").Value = 67
$ws.Range("AH22").Value = 11
$ws.Range("AI22").Value = 23
$ws.Range("AJ22").Value = 15
$ws.Range("AM22").Value = 41
$ws.Range("AO22").Value = 9.5
$ws.Range("AP22").Value = 21
$ws.Range("AQ22").Value = 34
$ws.Range("AR22").Value = 51
$ws.Range("AS22").Value = 151
$ws.Range("AT22").Value = 2.63
$ws.Range("AU22").Value = 8.5
$ws.Range("AV22").Value = 67
$ws.Range("AW22").Value = 6.5
$ws.Range("AX22").Value = 26
$ws.Range("AY22").Value = 34
$ws.Range("BA22").Value = 126
$ws.Range("BB22").Value = 301

$ws.Range("M28").Value = 1.04
$ws.Range("N28").Value = 13
$ws.Range("O28").Value = 1.25
$ws.Range("P28").Value = 3.75
$ws.Range("Q28").Value = 1.85
$ws.Range("R28").Value = 1.95

$ws.Range("G51").Value = 1.6
$ws.Range("I51").Value = 5.75
$ws.Range("J51").Value = 2.2
$ws.Range("L51").Value = 6
$ws.Range("N51").Value = 9.5
$ws.Range("S51").Value = 1.44
$ws.Range("T51").Value = 2.63
$ws.Range("U51").Value = 2.05
$ws.Range("V51").Value = 1.7
$ws.Range("Z51").Value = 11
$ws.Range("AE51").Value = 21
$ws.Range("AG51").Value = 501
$ws.Range("AK51").Value = 67
$ws.Range("AN51").Value = 3.4
$ws.Range("AQ51").Value = 26
$ws.Range("AT51").Value = 2.63
$ws.Range("AX51").Value = 34
$ws.Range("AZ51").Value = 126
$ws.Range("BB51").Value = 351

$ws.Range("G52").Value = 4.5
$ws.Range("H52").Value = 4.1
$ws.Range("I52").Value = 1.7
$ws.Range("L52").Value = 2.25
$ws.Range("M52").Value = 1.04
$ws.Range("N52").Value = 13
$ws.Range("Q52").Value = 1.75
$ws.Range("R52").Value = 2.05
$ws.Range("Y52").Value = 15
$ws.Range("AB52").Value = 41
$ws.Range("AR52").Value = 101
$ws.Range("AS52").Value = 201
$ws.Range("AU52").Value = 8
$ws.Range("AX52").Value = 8.5
$ws.Range("AZ52").Value = 26

$ws.Range("G78").Value = 2.7
$ws.Range("I78").Value = 2.8
$ws.Range("J78").Value = 3.5
$ws.Range("S78").Value = 1.57
$ws.Range("T78").Value = 2.25
$ws.Range("Z78").Value = 29
$ws.Range("AG78").Value = 1250
$ws.Range("AH78").Value = 7
$ws.Range("AL78").Value = 26
$ws.Range("AO78").Value = 17
$ws.Range("AS78").Value = 301
$ws.Range("AT78").Value = 2.25

$ws.Range("G79").Value = 2.55
$ws.Range("H79").Value = 3.05
$ws.Range("I79").Value = 2.72
$ws.Range("J79").Value = 3.1
$ws.Range("L79").Value = 3.25
$ws.Range("M79").Value = 1.05
$ws.Range("N79").Value = 9.800000000000001
$ws.Range("O79").Value = 1.31
$ws.Range("P79").Value = 2.9
$ws.Range("Q79").Value = 1.98
$ws.Range("R79").Value = 1.75
$ws.Range("U79").Value = 1.7
$ws.Range("V79").Value = 1.93
$ws.Range("W79").Value = 8.25
$ws.Range("X79").Value = 13
$ws.Range("Y79").Value = 9.5
$ws.Range("AA79").Value = 21
$ws.Range("AB79").Value = 29
$ws.Range("AC79").Value = 9
$ws.Range("AD79").Value = 5.9
$ws.Range("AF79").Value = 60
$ws.Range("AH79").Value = 8.5
$ws.Range("AI79").Value = 14
$ws.Range("AO79").Value = 13.5
$ws.Range("AP79").Value = 19.5
$ws.Range("AR79").Value = 80
$ws.Range("AS79").Value = 250
$ws.Range("AT79").Value = 2.52
$ws.Range("AX79").Value = 14.5

$ws.Range("I80").Value = 5.7
$ws.Range("J80").Value = 2.1
$ws.Range("K80").Value = 2.15
$ws.Range("U80").Value = 1.98
$ws.Range("V80").Value = 1.65
$ws.Range("W80").Value = 5.8
$ws.Range("X80").Value = 6.6
$ws.Range("AB80").Value = 32
$ws.Range("AE80").Value = 18.5
$ws.Range("AH80").Value = 14
$ws.Range("AI80").Value = 35
$ws.Range("AM80").Value = 65
$ws.Range("AO80").Value = 7.4
$ws.Range("AQ80").Value = 24
$ws.Range("AR80").Value = 60
$ws.Range("AU80").Value = 7.8
$ws.Range("AV80").Value = 75
$ws.Range("AW80").Value = 7.1
$ws.Range("AZ80").Value = 200

